# Refresh the "cryptos" price/volume snapshot (GitHub Actions bot update).
# Prices in column D that look like plain numbers are entered with a
# leading apostrophe so Excel keeps them as text (preserving trailing
# zeros such as "7.50" or "0.0930"), matching how the sheet originally
# stored every price as a text string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.782.10"
$ws.Range("E2").Value = "  +6.28%  "
$ws.Range("D3").Value = "2.399.38"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'113.97"
$ws.Range("E5").Value = "  +7.35%  "
$ws.Range("D6").Value = "'319.08"
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("D7").Value = "'0.635"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.626"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").Value = "'42.02"
$ws.Range("E10").Value = "  +5.61%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("E12").Value = "  +5.13%  "
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'15.85"
$ws.Range("E15").Value = "  +3.49%  "
$ws.Range("D16").Value = "2.763.63"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").Value = "2.394.70"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "45.698.25"
$ws.Range("E18").Value = "  +6.92%  "
$ws.Range("D19").Value = "'7.50"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").Value = "'0.0000109"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").Value = "'13.45"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "'74.80"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").Value = "'264.55"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  +5.28%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "'7.62"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").Value = "'11.32"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("E29").Value = "  +4.39%  "
$ws.Range("D30").Value = "'39.22"
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("D31").Value = "'22.79"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").Value = "'0.0972"
$ws.Range("E32").Value = "  +12.41%  "
$ws.Range("D33").Value = "'172.45"
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").Value = "'2.96"
$ws.Range("E34").Value = "  +5.09%  "
$ws.Range("D35").Value = "'0.133"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").Value = "'4.94"
$ws.Range("E36").Value = "  +6.57%  "
$ws.Range("E37").Value = "  +4.73%  "
$ws.Range("D38").Value = "'4.15"
$ws.Range("E38").Value = "  +13.65%  "
$ws.Range("E39").Value = "  +7.98%  "
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("E41").Value = "  +12.65%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.242"
$ws.Range("E42").Value = "  +5.54%  "
$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'100.65"
$ws.Range("E43").Value = "  -6.39%  "
$ws.Range("D44").Value = "'13.56"
$ws.Range("E44").Value = "  +9.73%  "
$ws.Range("D45").Value = "'72.19"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "'88.34"
$ws.Range("E46").Value = "  +14.95%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "'115.68"
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("D49").Value = "'5.73"
$ws.Range("E49").Value = "  +10.20%  "
$ws.Range("D50").Value = "'9.47"
$ws.Range("E50").Value = "  +6.49%  "
$ws.Range("D51").Value = "1.664.41"
$ws.Range("E51").Value = "  -3.37%  "
